$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Locate the paragraph that starts the "3.      Richter's Predictor..." item
# (item 3) and the paragraph that starts the "4.      Build a model..." item
# (item 4) so the deletion is robust even if paragraph indices shift.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $txt = $tr.Paragraphs($i).Text
    if ($txt.StartsWith("3.      Richter")) {
        $startPara = $i
    }
    if ($txt.StartsWith("4.      Build a model")) {
        $endPara = $i
    }
}

$deleteStart = $tr.Paragraphs($startPara).Start
$deleteLength = $tr.Paragraphs($endPara).Start - $deleteStart

# Remove the whole "3. Richter's Predictor: Modeling Earthquake Damage" item
# (its heading, its two description paragraphs and its link paragraph).
$rng = $tr.Characters($deleteStart, $deleteLength)
$rng.Delete()

# The former item "4." (Build a model that can describe what a video is
# about / YouTube-8m) now becomes the new item "3.".
$tr2 = $sh.TextFrame.TextRange
for ($i = 1; $i -le $tr2.Paragraphs().Count; $i++) {
    $para = $tr2.Paragraphs($i)
    if ($para.Text.StartsWith("4.      Build a model")) {
        $numberRun = $para.Runs(1)
        $numberRun.Text = "3.      "
        break
    }
}
